# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> (row -> new value) for column F updates
$sheetUpdates = @{
    "展览" = @{
        3  = 2173
        5  = 11468
        9  = 11413
        10 = 465
        11 = 1158
        12 = 74
        14 = 5678
    }
    "全部类型" = @{
        3  = 2173
        7  = 11468
        11 = 11413
        12 = 465
        13 = 1158
        14 = 74
        17 = 5678
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}

$wb.Save()
